$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "70.328.21"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.30%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.607.57"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.80%  "

$ws.Range("E4").Value = "  +0.08%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "580.36"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.12%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "189.68"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.63%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.603.88"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.77%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.629"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.50%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.189"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.74%  "

$ws.Range("E11").Value = "  -2.07%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "55.92"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -4.51%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000313"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +7.22%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.68"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.07%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.186.74"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.71%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "19.84"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.40%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.602.60"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.74%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "70.358.39"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.16%  "

$ws.Range("E19").Value = "  -0.52%  "

$ws.Range("E20").Value = "  +0.23%  "

$ws.Range("E21").Value = "  -2.17%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "494.29"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.93%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "19.19"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.58%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.95"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -7.01%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "96.23"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +4.92%  "

$ws.Range("E26").Value = "  -2.31%  "

$ws.Range("E27").Value = "  -5.75%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "11.10"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.82%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.47"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.50%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "32.21"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.50%  "

$ws.Range("E31").Value = "  -4.35%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "12.19"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.75%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.117"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.34%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "65.78"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.14%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "573.22"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -8.87%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "38.56"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -6.08%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0₃0816"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.41%  "

$ws.Range("E38").Value = "  +0.01%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.36"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +14.01%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.397"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.16%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.02"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.09%  "

$ws.Range("E42").Value = "  -1.67%  "

$ws.Range("E43").Value = "  -6.45%  "

$ws.Range("E44").Value = "  -4.47%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.226.04"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.44%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0443"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.52%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.78"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +5.73%  "

$ws.Range("E48").Value = "  +2.93%  "

$ws.Range("E49").Value = "  -0.42%  "

$ws.Range("E50").Value = "  +0.13%  "

$ws.Range("E51").Value = "  -3.64%  "
